$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/measure-report-assigned-practitioner"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet ---
$elements = $wb.Worksheets.Item("Elements")

# The Extension.url row's Fixed Value column mirrors the canonical URL, so
# it must be updated the same way as the Metadata sheet's URL cell.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/measure-report-assigned-practitioner"

# The "Constraint(s)" value on the top-level Extension row (row 2) was
# incorrectly duplicated from the Extension.extension row (row 4); the
# commit removes it from row 2 since it only applies to row 4.
$elements.Range("AI2").Value = ""
